$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unprotect the sheet so we can edit protected cells; re-protect at the end.
$ws.Unprotect()

# Update the confidentiality / as-of-date footer text (A80).
$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-27 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) figures for each holding row.
$ws.Range("D2").Value = 0.07487595161384469
$ws.Range("E2").Value = -0.01237682301931398
$ws.Range("D3").Value = 0.04588883799227831
$ws.Range("E3").Value = -0.01073454287079345
$ws.Range("D4").Value = 0.03887908173977962
$ws.Range("E4").Value = -0.008668336713189451
$ws.Range("D5").Value = 0.03510387327787091
$ws.Range("E5").Value = -0.0009248249438498224
$ws.Range("D6").Value = 0.03345308038852613
$ws.Range("E6").Value = -0.007406598300221612
$ws.Range("D7").Value = 0.03032496887254284
$ws.Range("E7").Value = 0.01557189643452994
$ws.Range("D8").Value = 0.03095922975366092
$ws.Range("E8").Value = 0.006856228030487799
$ws.Range("D9").Value = 0.02930553235330135
$ws.Range("E9").Value = -0.00153782456970486
$ws.Range("D10").Value = 0.02717610415688932
$ws.Range("E10").Value = -0.02273385461917177
$ws.Range("D11").Value = 0.02658695689040427
$ws.Range("E11").Value = 0.0009585430146179519
$ws.Range("D12").Value = 0.02314782846107857
$ws.Range("E12").Value = 0.0145633818779396
$ws.Range("D13").Value = 0.02402723944123579
$ws.Range("E13").Value = 0.01213130352045666
$ws.Range("D14").Value = 0.02057297944646204
$ws.Range("E14").Value = 0.0001890001890003656
$ws.Range("D15").Value = 0.02014400190835743
$ws.Range("E15").Value = -0.006976809086596147
$ws.Range("D16").Value = 0.0183827439065467
$ws.Range("E16").Value = 0.004811416921508771
$ws.Range("D17").Value = 0.01731660097630922
$ws.Range("E17").Value = -0.00337623971301948
$ws.Range("D18").Value = 0.01736850739896073
$ws.Range("E18").Value = -0.01207012811867847
$ws.Range("D19").Value = 0.01688485946804537
$ws.Range("E19").Value = 0.01553439540987611
$ws.Range("D20").Value = 0.01546249231592872
$ws.Range("E20").Value = -0.006447234475737917
$ws.Range("D21").Value = 0.01370484152941056
$ws.Range("E21").Value = 0.005414551607444862
$ws.Range("D22").Value = 0.01558270159437979
$ws.Range("E22").Value = -0.0207497820401048
$ws.Range("D23").Value = 0.01359123767117056
$ws.Range("E23").Value = -0.001792361781331797
$ws.Range("D24").Value = 0.0154861031796619
$ws.Range("E24").Value = -0.001234239248808056
$ws.Range("D25").Value = 0.01337934890794602
$ws.Range("E25").Value = 0.004537862792676428
$ws.Range("D26").Value = 0.01078754129675923
$ws.Range("E26").Value = 0.02054961089494167
$ws.Range("D27").Value = 0.0119560916637442
$ws.Range("E27").Value = 0.01200551689549267
$ws.Range("D28").Value = 0.01189828252515939
$ws.Range("E28").Value = -0.01073304407398956
$ws.Range("D29").Value = 0.01199296021485137
$ws.Range("E29").Value = 0.001195298492595809
$ws.Range("D30").Value = 0.01182979228155243
$ws.Range("E30").Value = 0.02708696340883887
$ws.Range("D31").Value = 0.01039946114762102
$ws.Range("E31").Value = 0.01423049894588879
$ws.Range("D32").Value = 0.01225094637869783
$ws.Range("E32").Value = 0.00861917326297279
$ws.Range("D33").Value = 0.01082520624604787
$ws.Range("E33").Value = -0.001242015613910574
$ws.Range("D34").Value = 0.01082754859364045
$ws.Range("E34").Value = 0.008359076867163306
$ws.Range("D35").Value = 0.01077573586489264
$ws.Range("E35").Value = 0.003391009477436535
$ws.Range("D36").Value = 0.009844886931602419
$ws.Range("E36").Value = -0.00907922912205561
$ws.Range("D37").Value = 0.0100472657636011
$ws.Range("E37").Value = 0.01126031612812994
$ws.Range("D38").Value = 0.008701305989954326
$ws.Range("E38").Value = 0.01892978857428984
$ws.Range("D39").Value = 0.01021038684994818
$ws.Range("E39").Value = -0.006937307297019468
$ws.Range("D40").Value = 0.00926763879088767
$ws.Range("E40").Value = -0.01066077602766047
$ws.Range("D41").Value = 0.008629021143247392
$ws.Range("E41").Value = 0.01561380268844048
$ws.Range("D42").Value = 0.008728805150691184
$ws.Range("E42").Value = 0.004132541888947294
$ws.Range("D43").Value = 0.009654032449759221
$ws.Range("E43").Value = -0.003532677264698525
$ws.Range("D44").Value = 0.008978030934541409
$ws.Range("E44").Value = -0.002410694718386841
$ws.Range("D45").Value = 0.008614030118654897
$ws.Range("E45").Value = 0.02782309817485706
$ws.Range("D46").Value = 0.009370139921538852
$ws.Range("E46").Value = 0.01295896328293722
$ws.Range("D47").Value = 0.008683176219587779
$ws.Range("E47").Value = -0.01696232034183609
$ws.Range("D48").Value = 0.008453907237226312
$ws.Range("E48").Value = 0.009431557481519359
$ws.Range("D49").Value = 0.007930017774670471
$ws.Range("E49").Value = 0.005848471422241985
$ws.Range("D50").Value = 0.008999486838489418
$ws.Range("E50").Value = -0.01030691708657805
$ws.Range("D51").Value = 0.007795192247241722
$ws.Range("E51").Value = 0.01501829960876688
$ws.Range("D52").Value = 0.008214425619361209
$ws.Range("E52").Value = -0.0003421806029223484
$ws.Range("D53").Value = 0.006674285230289778
$ws.Range("E53").Value = 0.01179195620130558
$ws.Range("D54").Value = 0.007532615082113807
$ws.Range("E54").Value = -0.008905915717200008
$ws.Range("D55").Value = 0.006625798635123428
$ws.Range("E55").Value = 0.01368119630925868
$ws.Range("D56").Value = 0.006800069296011179
$ws.Range("E56").Value = 0.001446729358272059
$ws.Range("D57").Value = 0.008013077420303263
$ws.Range("E57").Value = -0.004723820214208807
$ws.Range("D58").Value = 0.006554825503068334
$ws.Range("E58").Value = -0.003430531732418363
$ws.Range("D59").Value = 0.006536836273557341
$ws.Range("E59").Value = 0.01599587203302377
$ws.Range("D60").Value = 0.005846499591072958
$ws.Range("E60").Value = -0.01448717948717948
$ws.Range("D61").Value = 0.00571692092225158
$ws.Range("E61").Value = 0.02361636920858134
$ws.Range("D62").Value = 0.005793094065962194
$ws.Range("E62").Value = 0.01172569949862523
$ws.Range("D63").Value = 0.004910684880886472
$ws.Range("E63").Value = 0.0139662672670382
$ws.Range("D64").Value = 0.004956782281508393
$ws.Range("E64").Value = -0.006351126568879484
$ws.Range("D65").Value = 0.004535159714844478
$ws.Range("E65").Value = -0.004462441120568594
$ws.Range("D66").Value = 0.004474820840859687
$ws.Range("E66").Value = -0.001130653266331438
$ws.Range("D67").Value = 0.004518107424370516
$ws.Range("E67").Value = 0.002923976608186996
$ws.Range("D68").Value = 0.004367307086360388
$ws.Range("E68").Value = 0.01557522123893795
$ws.Range("D69").Value = 0.0041116632701065
$ws.Range("E69").Value = 0.001184942120135091
$ws.Range("D70").Value = 0.003516519593784459
$ws.Range("E70").Value = 0.004103165298944722
$ws.Range("D71").Value = 0.003673316341631584
$ws.Range("E71").Value = -0.007269388223591022
$ws.Range("D72").Value = 0.002941988576277097
$ws.Range("E72").Value = -0.01350318471337586
$ws.Range("D73").Value = 0.002353403473214272
$ws.Range("E73").Value = 0.002985906521219883
$ws.Range("D74").Value = 0.002375046764969686
$ws.Range("E74").Value = -0.0171604402540535
$ws.Range("D75").Value = 0.001926252966231874
$ws.Range("E75").Value = 0.02145045965270675
$ws.Range("D76").Value = 0.001974505326638966
$ws.Range("E76").Value = 0.01599126886210489
$ws.Range("E77").Value = 0.00004792443174417294

# Restore sheet protection.
$ws.Protect()
